$d = $word.ActiveDocument

# ------------------------------------------------------------------
# Edit 1: "Spring Boot" -> "Django"
# (the bold+italic occurrence in "...using Redux, React.js, and
# Spring Boot to monitor cluster statuses...")
# ------------------------------------------------------------------
$m1 = $d.Content
$m1.Find.Execute("Spring Boot") | Out-Null          # 1st occurrence (coursework list)
$m2 = $d.Range($m1.End, $d.Content.End)
$m2.Find.Execute("Spring Boot") | Out-Null          # 2nd occurrence (software tools list)
$m3 = $d.Range($m2.End, $d.Content.End)
$m3.Find.Execute("Spring Boot") | Out-Null          # 3rd occurrence - target (bold/italic)

$target1 = $d.Range($m3.Start, $m3.End)
$target1.Text = "Django"

# ------------------------------------------------------------------
# Edit 2: "Scrum Development, Circle CI/CD" ->
#         "Scrum Development, " + "Travis" + " CI/CD"  (Circle -> Travis,
#          split into three runs)
# ------------------------------------------------------------------

# "Circle" appears twice in the doc ("Circle/Travis CI (CI/CD)" in the
# Other Tools line, and inside "Scrum Development, Circle CI/CD").
# Locate the second occurrence specifically.
$c1 = $d.Content
$c1.Find.Execute("Circle") | Out-Null
$c2 = $d.Range($c1.End, $d.Content.End)
$c2.Find.Execute("Circle") | Out-Null

$circleRange = $d.Range($c2.Start, $c2.End)
$circleRange.Find.Execute("Circle", $true, $false, $false, $false, $false, $true, 0, $false, "Travis", 2) | Out-Null

# The phrase is now a single run: "Scrum Development, Travis CI/CD".
# Split it into three runs by re-toggling Bold on just the "Travis" word
# (forces Word to break the run at the word boundaries while keeping
# identical formatting on every piece).
$phrase = $d.Content
$phrase.Find.Execute("Scrum Development, Travis CI/CD") | Out-Null

$travisWord = $d.Range($phrase.Start, $phrase.End)
$travisWord.Find.Execute("Travis", $true, $false, $false, $false, $false, $true, 0, $false) | Out-Null
$tStart = $travisWord.Start
$tEnd = $travisWord.End

$toggleOff = $d.Range($tStart, $tEnd)
$toggleOff.Font.Bold = $false
$toggleOn = $d.Range($tStart, $tEnd)
$toggleOn.Font.Bold = $true

Write-Output "done"
